# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text assignment (safe for non-numeric-looking strings: names, URLs,
# and percentage strings that include "%"/spaces so Excel won't coerce them
# into numeric values).
function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Text assignment that forces the cell to be stored as text even when the
# content looks like a number (e.g. "230.00", "1.20", "0.0240"). Without
# this, Excel's COM layer silently reinterprets such strings as doubles,
# which both normalizes away formatting (trailing zeros) and can introduce
# binary floating point artifacts (e.g. "57.19" -> "57.189999999999998").
# Setting the NumberFormat to Text ("@") first preserves the literal
# characters; restoring the style to "Normal" afterwards avoids leaving a
# stray number format applied to the cell (matching the original,
# style-less inline-string cells).
function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "41.712.29"
Set-Cell "E2" "  +0.09%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.210.01"
Set-Cell "E3" "  -2.34%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  -0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "230.00"
Set-Cell "E5" "  -2.60%  "

# Row 6 - XRP
Set-Cell "E6" "  -3.80%  "

# Row 7 - Solana
Set-TextCell "D7" "60.35"
Set-Cell "E7" "  -5.85%  "

# Row 8 - USDC
Set-Cell "E8" "  -0.06%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.403"
Set-Cell "E9" "  -2.30%  "

# Row 10 - OKB
Set-TextCell "D10" "57.19"
Set-Cell "E10" "  -4.37%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.0889"
Set-Cell "E11" "  -1.51%  "

# Row 12 - TRON
Set-Cell "E12" "  -2.40%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "2.534.83"
Set-Cell "E13" "  -2.47%  "

# Row 14 - Chainlink
Set-Cell "E14" "  -4.49%  "

# Row 15 - Avalanche
Set-TextCell "D15" "22.19"
Set-Cell "E15" "  -3.06%  "

# Row 16 - Polkadot
Set-Cell "E16" "  -0.66%  "

# Row 17 - Polygon
Set-Cell "E17" "  -3.98%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "2.217.91"
Set-Cell "E18" "  -1.42%  "

# Row 19 - WrappedBTC
Set-TextCell "D19" "41.750.72"
Set-Cell "E19" "  +0.34%  "

# Row 20 - Litecoin
Set-TextCell "D20" "72.18"
Set-Cell "E20" "  -3.64%  "

# Row 21 - ShibaInu (contains subscript-3 U+2083 character)
$subscriptThree = [char]0x2083
$d21Value = [string]::Concat("0.0", $subscriptThree, "0902")
Set-TextCell "D21" $d21Value
Set-Cell "E21" "  -4.09%  "

# Row 22 - Uniswap
Set-TextCell "D22" "6.05"
Set-Cell "E22" "  -2.39%  "

# Row 23 - BitcoinCash
Set-TextCell "D23" "242.63"
Set-Cell "E23" "  -3.86%  "

# Row 25/26 - PancakeSwap and Toncoin swap positions
Set-Cell "B25" "Toncoin"
Set-Cell "C25" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D25" "2.47"
Set-Cell "E25" "  +4.83%  "

Set-Cell "B26" "PancakeSwap"
Set-Cell "C26" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D26" "2.35"
Set-Cell "E26" "  -3.93%  "

# Row 27 - Cosmos
Set-Cell "E27" "  -2.15%  "

# Row 28 - Monero
Set-TextCell "D28" "169.33"
Set-Cell "E28" "  -1.36%  "

# Row 29 - Kaspa
Set-Cell "E29" "  -5.91%  "

# Row 30 - ImmutableX
Set-TextCell "D30" "1.45"
Set-Cell "E30" "  -0.40%  "

# Row 31 - EthereumClassic
Set-TextCell "D31" "19.76"
Set-Cell "E31" "  -3.93%  "

# Row 32 - WEMIXToken
Set-TextCell "D32" "2.61"
Set-Cell "E32" "  -7.88%  "

# Row 33 - Stellar
Set-Cell "E33" "  -3.97%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextCell "D34" "5.02"
Set-Cell "E34" "  -2.06%  "

# Row 35 - Filecoin
Set-Cell "E35" "  -4.06%  "

# Row 36 - Hedera
Set-TextCell "D36" "0.0649"

# Row 37 - LidoDAOToken
Set-TextCell "D37" "2.37"
Set-Cell "E37" "  -4.08%  "

# Row 38 - THORChain
Set-TextCell "D38" "6.33"
Set-Cell "E38" "  -8.55%  "

# Row 39 - RenderToken
Set-Cell "E39" "  -8.44%  "

# Row 40 - TerraClassic
Set-TextCell "D40" "0.000239"
Set-Cell "E40" "  -10.56%  "

# Row 41 - BinanceUSD
Set-TextCell "D41" "0.999"
Set-Cell "E41" "  -0.04%  "

# Row 42 - VeChain
Set-TextCell "D42" "0.0240"
Set-Cell "E42" "  -0.90%  "

# Row 43 - FraxShare
Set-TextCell "D43" "8.59"
Set-Cell "E43" "  -1.62%  "

# Row 44 - Cronos
Set-TextCell "D44" "0.0957"
Set-Cell "E44" "  -3.23%  "

# Row 45/47 - TrustWalletToken and FTXToken swap positions
Set-Cell "B45" "FTXToken"
Set-Cell "C45" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D45" "4.42"
Set-Cell "E45" "  -13.25%  "

# Row 46 - Aave
Set-TextCell "D46" "97.34"
Set-Cell "E46" "  -5.49%  "

Set-Cell "B47" "TrustWalletToken"
Set-Cell "C47" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D47" "1.20"
Set-Cell "E47" "  -3.67%  "

# Row 48 - Maker
Set-TextCell "D48" "1.466.49"
Set-Cell "E48" "  -2.99%  "

# Row 49 - InjectiveProtocol
Set-TextCell "D49" "16.40"
Set-Cell "E49" "  -7.01%  "

# Row 50 - HuobiToken
Set-Cell "E50" "  -1.43%  "

# Row 51 - ARBITRUM
Set-Cell "E51" "  -5.44%  "
